$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text "Ready for handoff" -> "In Translation" on all sheets
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# Shrink column widths for zh-cn/de-de status columns (Overview E:F, zh-cn/de-de C)
# Target stored width ~13.4101845877511 chars; use 12.5 as the ColumnWidth input,
# which this engine's rounding resolves to the closest representable stored width.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
